$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.817.91"
$ws.Range("E2").Value = "  +1.22%  "
$ws.Range("D3").Value = "2.540.69"
$ws.Range("E3").Value = "  +0.54%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'591.32"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.39%  "
$ws.Range("D6").Value = "'173.38"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.49%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("E8").Value = "  -0.15%  "
$ws.Range("D9").Value = "2.540.42"
$ws.Range("E9").Value = "  +0.56%  "
$ws.Range("E10").Value = "  +0.36%  "
$ws.Range("E12").Value = "  -0.24%  "
$ws.Range("E13").Value = "  +0.49%  "
$ws.Range("D14").Value = "'26.55"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.11%  "
$ws.Range("E15").Value = "  +1.13%  "
$ws.Range("D16").Value = "2.919.75"
$ws.Range("D17").Value = "67.480.18"
$ws.Range("E17").Value = "  +0.93%  "
$ws.Range("D18").Value = "2.533.61"
$ws.Range("E18").Value = "  +0.01%  "
$ws.Range("D19").Value = "'11.81"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.11%  "
$ws.Range("D20").Value = "'7.95"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.81%  "
$ws.Range("D21").Value = "'370.39"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.33%  "
$ws.Range("B22").Value = "Polkadot"
$ws.Range("C22").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D22").Value = "'4.15"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.46%  "
$ws.Range("B23").Value = "NEARProtocol"
$ws.Range("C23").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D23").Value = "'4.59"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.42%  "
$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").Value = "'71.67"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.75%  "
$ws.Range("D25").Value = "'1.00"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.02%  "
$ws.Range("B26").Value = "SuiNetwork"
$ws.Range("C26").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D26").Value = "'1.93"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.71%  "
$ws.Range("B27").Value = "Aptos"
$ws.Range("C27").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D27").Value = "'10.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.26%  "
$ws.Range("B28").Value = "WrappedeETH"
$ws.Range("C28").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D28").Value = "2.666.30"
$ws.Range("E28").Value = "  +0.50%  "
$ws.Range("B29").Value = "PEPE"
$ws.Range("C29").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D29").Value = "0.0₃0970"
$ws.Range("E29").Value = "  -0.34%  "
$ws.Range("B30").Value = "InternetComputer(DFINITY)"
$ws.Range("C30").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D30").Value = "'8.50"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.85%  "
$ws.Range("B31").Value = "Bittensor"
$ws.Range("C31").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D31").Value = "'538.45"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.27%  "
$ws.Range("B32").Value = "Fetch.AI"
$ws.Range("C32").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D32").Value = "'1.32"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.24%  "
$ws.Range("B33").Value = "PancakeSwap"
$ws.Range("C33").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D33").Value = "'1.88"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.75%  "
$ws.Range("B34").Value = "Kaspa"
$ws.Range("C34").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D34").Value = "'0.129"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.56%  "
$ws.Range("B35").Value = "FirstDigitalUSD"
$ws.Range("C35").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D35").Value = "'0.999"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.05%  "
$ws.Range("B36").Value = "Monero"
$ws.Range("C36").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D36").Value = "'159.33"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.36%  "
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").Value = "'1.44"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.11%  "
$ws.Range("B38").Value = "EthereumClassic"
$ws.Range("C38").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D38").Value = "'19.21"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.37%  "
$ws.Range("B39").Value = "WhiteBITCoin"
$ws.Range("C39").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D39").Value = "'18.63"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.99%  "
$ws.Range("B40").Value = "PolygonEcosystemToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D40").Value = "'0.352"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.51%  "
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D41").Value = "'5.17"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.79%  "
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").Value = "'1.78"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.51%  "
$ws.Range("B43").Value = "dogwifhat"
$ws.Range("C43").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D43").Value = "'2.58"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.70%  "
$ws.Range("B44").Value = "USDe"
$ws.Range("C44").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D44").Value = "'1.00"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.05%  "
$ws.Range("B45").Value = "BabyDogeCoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D45").Value = "0.0₆0290"
$ws.Range("E45").Value = "  +4.66%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "'147.52"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.04%  "
$ws.Range("B47").Value = "Filecoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D47").Value = "'3.71"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.92%  "
$ws.Range("D48").Value = "'0.553"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.28%  "
$ws.Range("B49").Value = "Optimism"
$ws.Range("C49").Value = "https://coinranking.com/coin/n1p-s_gm1+optimism-op"
$ws.Range("D49").Value = "'1.72"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.23%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "'0.0747"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.38%  "
$ws.Range("B51").Value = "Mantle"
$ws.Range("C51").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D51").Value = "'0.598"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.45%  "
